# Outline-of-thesis diagram: reword the title textbox so "sagittal-plane"
# becomes "sagittal plane" (hyphen -> space), matching the target commit.
#
# The shape is slide 1, "TextBox 76" (shape index 6), whose first paragraph
# currently reads (as a single run):
#   "Development of a model for the average sagittal-plane hip and knee
#    angle functions from the RISC data"
#
# Retyping just the "average sagittal-plane " portion of that paragraph
# naturally splits it into the three runs seen in the diff:
#   1) "Development of a model for the "
#   2) "average sagittal plane "
#   3) "hip and knee angle functions from the RISC data"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $null
try {
    $sh = $s.Shapes.Item("TextBox 76")
} catch {
    $sh = $null
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(6)
}

$tf = $sh.TextFrame
$tr = $tf.TextRange

$titlePara = $tr.Paragraphs(1, 1)

$oldFragment = "average sagittal-plane "
$newFragment = "average sagittal plane "

$idx = $titlePara.Text.IndexOf($oldFragment)
if ($idx -ge 0) {
    $target = $titlePara.Characters($idx + 1, $oldFragment.Length)
    $target.Text = $newFragment
}
